$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A dates: rows 2-10 -> 12.09.11, rows 11-18 -> 12.010.11
# Force text format (so Excel doesn't auto-convert these date-like strings
# into date serials), set the value, then restore the cell's original
# (unstyled) appearance by copying the style from an untouched neighbor cell.
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = "12.09.11"
    $cell.Style = $ws.Cells.Item($r, 3).Style
}
for ($r = 11; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = "12.010.11"
    $cell.Style = $ws.Cells.Item($r, 3).Style
}

# Update selection to A12:A18 with active cell A12
$ws.Range("A12:A18").Select()
